# Auto-generated script applying the diff changes to all 8 sheets
$wb = $excel.ActiveWorkbook

# ===== Sheet 1: ALC =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("H5").Value = 299.5
$ws.Range("I5").Value = 121.2
$ws.Range("J5").Value = 1191
$ws.Range("K5").Value = 121.2
$ws.Range("L5").Value = 1191
$ws.Range("M5").Value = -6.200000000000003
$ws.Range("N5").Value = -1421
$ws.Range("H12").Value = 3788420
$ws.Range("J12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("N12").Value = -2340
$ws.Range("H17").Value = 332012.94
$ws.Range("J17").Value = 375734.66
$ws.Range("L17").Value = 1127203.98
$ws.Range("N17").Value = -1127539.98
$ws.Range("H80").Value = 62500932
$ws.Range("I80").Value = 76923800
$ws.Range("K80").Value = 230771400
$ws.Range("M80").Value = -230770402
$ws.Range("H82").Value = 2833.3635
$ws.Range("I82").Value = 2456.7
$ws.Range("J82").Value = 6600
$ws.Range("K82").Value = 7370.099999999999
$ws.Range("L82").Value = 19800
$ws.Range("M82").Value = -6964.099999999999
$ws.Range("N82").Value = -20612
$ws.Range("H83").Value = 62500932
$ws.Range("I83").Value = 76923800
$ws.Range("K83").Value = 692314200
$ws.Range("M83").Value = -692309208
$ws.Range("H85").Value = 2833.3635
$ws.Range("I85").Value = 2456.7
$ws.Range("J85").Value = 6600
$ws.Range("K85").Value = 7370.099999999999
$ws.Range("L85").Value = 19800
$ws.Range("M85").Value = -5966.099999999999
$ws.Range("N85").Value = -22608
$ws.Range("H88").Value = 960
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 933.3333
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 933.3333
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -1745.3333
$ws.Range("H91").Value = 960
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 933.3333
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 933.3333
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -3741.3333
$ws.Range("H98").Value = 848.70966
$ws.Range("I98").Value = 870.6539
$ws.Range("J98").Value = 734.6
$ws.Range("K98").Value = 870.6539
$ws.Range("L98").Value = 734.6
$ws.Range("M98").Value = 627.3461
$ws.Range("N98").Value = -3730.6
$ws.Range("H99").Value = 625.1667
$ws.Range("I99").Value = 677.9231
$ws.Range("K99").Value = 2033.7693
$ws.Range("M99").Value = -535.7692999999999
$ws.Range("H101").Value = 233.1
$ws.Range("I101").Value = 237.33333
$ws.Range("K101").Value = 711.99999
$ws.Range("M101").Value = 910.00001
$ws.Range("H111").Value = 3620.375
$ws.Range("I111").Value = 3620.375
$ws.Range("K111").Value = 10861.125
$ws.Range("M111").Value = -7794.125
$ws.Range("H112").Value = 4188.364
$ws.Range("I112").Value = 2036.6
$ws.Range("J112").Value = 4572.607
$ws.Range("K112").Value = 6109.799999999999
$ws.Range("L112").Value = 13717.821
$ws.Range("M112").Value = -5001.799999999999
$ws.Range("N112").Value = -15933.821
$ws.Range("H113").Value = 5396.7896
$ws.Range("I113").Value = 4064.4443
$ws.Range("K113").Value = 4064.4443
$ws.Range("M113").Value = -810.4443000000001
$ws.Range("H116").Value = 6968.6875
$ws.Range("I116").Value = 6359.9
$ws.Range("J116").Value = 7983.3335
$ws.Range("K116").Value = 6359.9
$ws.Range("L116").Value = 7983.3335
$ws.Range("M116").Value = -2917.9
$ws.Range("N116").Value = -14867.3335
$ws.Range("H122").Value = 848.70966
$ws.Range("I122").Value = 870.6539
$ws.Range("J122").Value = 734.6
$ws.Range("K122").Value = 2611.9617
$ws.Range("L122").Value = 2203.8
$ws.Range("M122").Value = -161.9616999999998
$ws.Range("N122").Value = -7103.8
$ws.Range("H128").Value = 115892
$ws.Range("J128").Value = 115892
$ws.Range("L128").Value = 115892
$ws.Range("N128").Value = -125852
$ws.Range("H132").Value = 21652.062
$ws.Range("I132").Value = 23603.613
$ws.Range("J132").Value = 185
$ws.Range("K132").Value = 70810.83900000001
$ws.Range("L132").Value = 555
$ws.Range("M132").Value = -68280.83900000001
$ws.Range("N132").Value = -5615
$ws.Range("H135").Value = 1441.1578
$ws.Range("I135").Value = 786.4666999999999
$ws.Range("J135").Value = 3896.25
$ws.Range("K135").Value = 7078.2003
$ws.Range("L135").Value = 35066.25
$ws.Range("M135").Value = -4543.2003
$ws.Range("N135").Value = -40136.25
$ws.Range("H137").Value = 5559100
$ws.Range("I137").Value = 1834.0834
$ws.Range("K137").Value = 5502.2502
$ws.Range("M137").Value = -2952.2502
$ws.Range("H138").Value = 3619.9888
$ws.Range("I138").Value = 4011.4807
$ws.Range("J138").Value = 3084.2632
$ws.Range("K138").Value = 12034.4421
$ws.Range("L138").Value = 9252.7896
$ws.Range("M138").Value = -6894.4421
$ws.Range("N138").Value = -19532.7896
$ws.Range("H141").Value = 1348.4348
$ws.Range("I141").Value = 1487.9
$ws.Range("J141").Value = 418.66666
$ws.Range("K141").Value = 4463.700000000001
$ws.Range("L141").Value = 1255.99998
$ws.Range("M141").Value = 716.2999999999993
$ws.Range("N141").Value = -11615.99998

# ===== Sheet 2: ARM =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 2641.4
$ws.Range("I2").Value = 934.75
$ws.Range("J2").Value = 3779.1667
$ws.Range("K2").Value = 934.75
$ws.Range("L2").Value = 3779.1667
$ws.Range("M2").Value = -821.75
$ws.Range("N2").Value = -4005.1667
$ws.Range("H5").Value = 783.5
$ws.Range("I5").Value = 628
$ws.Range("K5").Value = 628
$ws.Range("M5").Value = -516
$ws.Range("H8").Value = 28900
$ws.Range("H13").Value = 3500
$ws.Range("J13").Value = 3200
$ws.Range("L13").Value = 3200
$ws.Range("N13").Value = -3488
$ws.Range("H32").Value = 1131.7024
$ws.Range("I32").Value = 944.8875
$ws.Range("J32").Value = 4868
$ws.Range("K32").Value = 944.8875
$ws.Range("L32").Value = 4868
$ws.Range("M32").Value = -657.8875
$ws.Range("N32").Value = -5442
$ws.Range("H61").Value = 1552805.2
$ws.Range("I61").Value = 3987.9285
$ws.Range("J61").Value = 5889494
$ws.Range("K61").Value = 3987.9285
$ws.Range("L61").Value = 5889494
$ws.Range("M61").Value = -3775.9285
$ws.Range("N61").Value = -5889918
$ws.Range("H63").Value = 16622.346
$ws.Range("I63").Value = 2700.3635
$ws.Range("K63").Value = 2700.3635
$ws.Range("M63").Value = -2014.3635
$ws.Range("H66").Value = 16622.346
$ws.Range("I66").Value = 2700.3635
$ws.Range("K66").Value = 13501.8175
$ws.Range("M66").Value = -10069.8175
$ws.Range("H74").Value = 544400.8
$ws.Range("I74").Value = 1341.1562
$ws.Range("K74").Value = 1341.1562
$ws.Range("M74").Value = -467.1561999999999
$ws.Range("H77").Value = 544400.8
$ws.Range("I77").Value = 1341.1562
$ws.Range("K77").Value = 6705.780999999999
$ws.Range("M77").Value = -2337.780999999999
$ws.Range("H88").Value = 3742.1667
$ws.Range("I88").Value = 3151
$ws.Range("K88").Value = 3151
$ws.Range("M88").Value = -2745
$ws.Range("H91").Value = 3742.1667
$ws.Range("I91").Value = 3151
$ws.Range("K91").Value = 3151
$ws.Range("M91").Value = -1747
$ws.Range("H96").Value = 6886.6
$ws.Range("J96").Value = 6886.6
$ws.Range("L96").Value = 6886.6
$ws.Range("N96").Value = -12378.6
$ws.Range("H109").Value = 63762.547
$ws.Range("J109").Value = 63762.547
$ws.Range("L109").Value = 63762.547
$ws.Range("N109").Value = -66536.54699999999
$ws.Range("H116").Value = 2641.4
$ws.Range("I116").Value = 934.75
$ws.Range("J116").Value = 3779.1667
$ws.Range("K116").Value = 934.75
$ws.Range("L116").Value = 3779.1667
$ws.Range("M116").Value = 1359.25
$ws.Range("N116").Value = -8367.1667
$ws.Range("H122").Value = 2760.1738
$ws.Range("I122").Value = 1751.9474
$ws.Range("K122").Value = 5255.8422
$ws.Range("M122").Value = -2805.8422
$ws.Range("H132").Value = 1213.7368
$ws.Range("I132").Value = 1238.8823
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3716.6469
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1186.6469
$ws.Range("N132").Value = -8060
$ws.Range("H136").Value = 1552805.2
$ws.Range("I136").Value = 3987.9285
$ws.Range("J136").Value = 5889494
$ws.Range("K136").Value = 11963.7855
$ws.Range("L136").Value = 17668482
$ws.Range("M136").Value = -9413.7855
$ws.Range("N136").Value = -17673582
$ws.Range("H141").Value = 111384.664
$ws.Range("J141").Value = 111384.664
$ws.Range("L141").Value = 111384.664
$ws.Range("N141").Value = -121744.664

# ===== Sheet 3: BSM =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 2641.4
$ws.Range("I3").Value = 934.75
$ws.Range("J3").Value = 3779.1667
$ws.Range("K3").Value = 934.75
$ws.Range("L3").Value = 3779.1667
$ws.Range("M3").Value = -820.75
$ws.Range("N3").Value = -4007.1667
$ws.Range("H4").Value = 783.5
$ws.Range("I4").Value = 628
$ws.Range("K4").Value = 628
$ws.Range("M4").Value = -513
$ws.Range("H5").Value = 862.3333
$ws.Range("I5").Value = 293.5
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 293.5
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -180.5
$ws.Range("N5").Value = -2226
$ws.Range("H7").Value = 10000617
$ws.Range("I7").Value = 15000250
$ws.Range("J7").Value = 1350
$ws.Range("K7").Value = 15000250
$ws.Range("L7").Value = 1350
$ws.Range("M7").Value = -15000137
$ws.Range("N7").Value = -1576
$ws.Range("H20").Value = 1612.6
$ws.Range("I20").Value = 1600
$ws.Range("J20").Value = 1618
$ws.Range("K20").Value = 1600
$ws.Range("L20").Value = 1618
$ws.Range("M20").Value = -1353
$ws.Range("N20").Value = -2112
$ws.Range("H35").Value = 68996.336
$ws.Range("J35").Value = 68996.336
$ws.Range("L35").Value = 68996.336
$ws.Range("N35").Value = -69616.336
$ws.Range("H82").Value = 25970.182
$ws.Range("I82").Value = 11546.75
$ws.Range("J82").Value = 64432.668
$ws.Range("K82").Value = 11546.75
$ws.Range("L82").Value = 64432.668
$ws.Range("M82").Value = -11163.75
$ws.Range("N82").Value = -65198.668
$ws.Range("H85").Value = 25970.182
$ws.Range("I85").Value = 11546.75
$ws.Range("J85").Value = 64432.668
$ws.Range("K85").Value = 11546.75
$ws.Range("L85").Value = 64432.668
$ws.Range("M85").Value = -10220.75
$ws.Range("N85").Value = -67084.66800000001
$ws.Range("H86").Value = 5450
$ws.Range("I86").Value = 1933.3334
$ws.Range("J86").Value = 16000
$ws.Range("K86").Value = 1933.3334
$ws.Range("L86").Value = 16000
$ws.Range("M86").Value = -810.3334
$ws.Range("N86").Value = -18246
$ws.Range("H89").Value = 5450
$ws.Range("I89").Value = 1933.3334
$ws.Range("J89").Value = 16000
$ws.Range("K89").Value = 9666.666999999999
$ws.Range("L89").Value = 80000
$ws.Range("M89").Value = -4050.666999999999
$ws.Range("N89").Value = -91232
$ws.Range("H99").Value = 5587.5312
$ws.Range("I99").Value = 8070.8887
$ws.Range("K99").Value = 8070.8887
$ws.Range("M99").Value = -6572.8887
$ws.Range("H124").Value = 88994
$ws.Range("J124").Value = 88994
$ws.Range("L124").Value = 88994
$ws.Range("N124").Value = -98814
$ws.Range("H134").Value = 22502264
$ws.Range("I134").Value = 1829.84
$ws.Range("J134").Value = 60002988
$ws.Range("K134").Value = 5489.52
$ws.Range("L134").Value = 180008964
$ws.Range("M134").Value = -2954.52
$ws.Range("N134").Value = -180014034

# ===== Sheet 4: CRP =====
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 33.6
$ws.Range("I7").Value = 27
$ws.Range("K7").Value = 27
$ws.Range("M7").Value = 86
$ws.Range("H9").Value = 34965
$ws.Range("J9").Value = 34965
$ws.Range("L9").Value = 34965
$ws.Range("N9").Value = -35301
$ws.Range("H31").Value = 4467.5386
$ws.Range("I31").Value = 3775.7693
$ws.Range("J31").Value = 4813.423
$ws.Range("K31").Value = 3775.7693
$ws.Range("L31").Value = 4813.423
$ws.Range("M31").Value = -3480.7693
$ws.Range("N31").Value = -5403.423
$ws.Range("H34").Value = 4467.5386
$ws.Range("I34").Value = 3775.7693
$ws.Range("J34").Value = 4813.423
$ws.Range("K34").Value = 3775.7693
$ws.Range("L34").Value = 4813.423
$ws.Range("M34").Value = -3573.7693
$ws.Range("N34").Value = -5217.423
$ws.Range("H58").Value = 2387.6667
$ws.Range("I58").Value = 2118.6155
$ws.Range("J58").Value = 3087.2
$ws.Range("K58").Value = 2118.6155
$ws.Range("L58").Value = 3087.2
$ws.Range("M58").Value = -1915.6155
$ws.Range("N58").Value = -3493.2
$ws.Range("H99").Value = 1822191.4
$ws.Range("I99").Value = 2860567.8
$ws.Range("J99").Value = 5032.5
$ws.Range("K99").Value = 2860567.8
$ws.Range("L99").Value = 5032.5
$ws.Range("M99").Value = -2859069.8
$ws.Range("N99").Value = -8028.5
$ws.Range("H105").Value = 2104.5715
$ws.Range("J105").Value = 1976.6154
$ws.Range("L105").Value = 1976.6154
$ws.Range("N105").Value = -5470.6154
$ws.Range("H107").Value = 2606.1667
$ws.Range("I107").Value = 2111.125
$ws.Range("J107").Value = 3596.25
$ws.Range("K107").Value = 2111.125
$ws.Range("L107").Value = 3596.25
$ws.Range("M107").Value = -191.125
$ws.Range("N107").Value = -7436.25
$ws.Range("H126").Value = 1822191.4
$ws.Range("I126").Value = 2860567.8
$ws.Range("J126").Value = 5032.5
$ws.Range("K126").Value = 8581703.399999999
$ws.Range("L126").Value = 15097.5
$ws.Range("M126").Value = -8579233.399999999
$ws.Range("N126").Value = -20037.5
$ws.Range("H132").Value = 2856.9062
$ws.Range("I132").Value = 2558.92
$ws.Range("J132").Value = 3921.1428
$ws.Range("K132").Value = 7676.76
$ws.Range("L132").Value = 11763.4284
$ws.Range("M132").Value = -5146.76
$ws.Range("N132").Value = -16823.4284
$ws.Range("H134").Value = 2096.182
$ws.Range("I134").Value = 1561.5555
$ws.Range("J134").Value = 2466.3076
$ws.Range("K134").Value = 4684.666499999999
$ws.Range("L134").Value = 7398.9228
$ws.Range("M134").Value = -2149.666499999999
$ws.Range("N134").Value = -12468.9228
$ws.Range("H136").Value = 2387.6667
$ws.Range("I136").Value = 2118.6155
$ws.Range("J136").Value = 3087.2
$ws.Range("K136").Value = 6355.8465
$ws.Range("L136").Value = 9261.599999999999
$ws.Range("M136").Value = -3805.8465
$ws.Range("N136").Value = -14361.6

# ===== Sheet 5: CUL =====
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 41.05
$ws.Range("I2").Value = 41.72222
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 250.33332
$ws.Range("L2").Value = 210
$ws.Range("M2").Value = -137.33332
$ws.Range("N2").Value = -436
$ws.Range("H4").Value = 6345586
$ws.Range("I4").Value = 6769768.5
$ws.Range("J4").Value = 5126062.5
$ws.Range("K4").Value = 20309305.5
$ws.Range("L4").Value = 15378187.5
$ws.Range("M4").Value = -20309193.5
$ws.Range("N4").Value = -15378411.5
$ws.Range("H11").Value = 90494
$ws.Range("I11").Value = 90494
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 271482
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -271342
$ws.Range("N11").ClearContents()
$ws.Range("H12").Value = 1102.6666
$ws.Range("J12").Value = 1468.1428
$ws.Range("L12").Value = 4404.428400000001
$ws.Range("N12").Value = -4750.428400000001
$ws.Range("H32").Value = 1665
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 9000
$ws.Range("N32").Value = -9566
$ws.Range("H68").Value = 2174.9
$ws.Range("I68").Value = 2399.6667
$ws.Range("J68").Value = 2078.5715
$ws.Range("K68").Value = 7199.000100000001
$ws.Range("L68").Value = 6235.7145
$ws.Range("M68").Value = -6388.000100000001
$ws.Range("N68").Value = -7857.7145
$ws.Range("H70").Value = 7615.8335
$ws.Range("I70").Value = 3975
$ws.Range("K70").Value = 11925
$ws.Range("M70").Value = -11610
$ws.Range("H71").Value = 2174.9
$ws.Range("I71").Value = 2399.6667
$ws.Range("J71").Value = 2078.5715
$ws.Range("K71").Value = 21597.0003
$ws.Range("L71").Value = 18707.1435
$ws.Range("M71").Value = -17541.0003
$ws.Range("N71").Value = -26819.1435
$ws.Range("H73").Value = 7615.8335
$ws.Range("I73").Value = 3975
$ws.Range("K73").Value = 11925
$ws.Range("M73").Value = -10833
$ws.Range("H80").Value = 1499.5
$ws.Range("I80").Value = 1499.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4498.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3562.5
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 15878826
$ws.Range("I81").Value = 37039760
$ws.Range("J81").Value = 8127.4165
$ws.Range("K81").Value = 111119280
$ws.Range("L81").Value = 24382.2495
$ws.Range("M81").Value = -111118157
$ws.Range("N81").Value = -26628.2495
$ws.Range("H83").Value = 1499.5
$ws.Range("I83").Value = 1499.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13495.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8815.5
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 15878826
$ws.Range("I84").Value = 37039760
$ws.Range("J84").Value = 8127.4165
$ws.Range("K84").Value = 333357840
$ws.Range("L84").Value = 73146.7485
$ws.Range("M84").Value = -333352224
$ws.Range("N84").Value = -84378.7485
$ws.Range("H113").Value = 1582.2307
$ws.Range("I113").Value = 1130
$ws.Range("J113").Value = 1717.9
$ws.Range("K113").Value = 3390
$ws.Range("L113").Value = 5153.700000000001
$ws.Range("M113").Value = -1220
$ws.Range("N113").Value = -9493.700000000001
$ws.Range("H114").Value = 1486.5714
$ws.Range("J114").Value = 1566.875
$ws.Range("L114").Value = 4700.625
$ws.Range("N114").Value = -11208.625
$ws.Range("H117").Value = 15154120
$ws.Range("J117").Value = 19610382
$ws.Range("L117").Value = 58831146
$ws.Range("N117").Value = -58838030
$ws.Range("H119").Value = 100023816
$ws.Range("I119").Value = 166677200
$ws.Range("K119").Value = 500031600
$ws.Range("M119").Value = -500026762
$ws.Range("H121").Value = 7027.316
$ws.Range("I121").Value = 842.4286
$ws.Range("J121").Value = 10635.167
$ws.Range("K121").Value = 2527.2858
$ws.Range("L121").Value = 31905.501
$ws.Range("M121").Value = -1217.2858
$ws.Range("N121").Value = -34525.501
$ws.Range("H129").Value = 3005.7727
$ws.Range("I129").Value = 621.5714
$ws.Range("J129").Value = 7178.125
$ws.Range("K129").Value = 1864.7142
$ws.Range("L129").Value = 21534.375
$ws.Range("M129").Value = 3135.2858
$ws.Range("N129").Value = -31534.375
$ws.Range("H131").Value = 2675918.8
$ws.Range("I131").Value = 5349091
$ws.Range("J131").Value = 2746.5881
$ws.Range("K131").Value = 16047273
$ws.Range("L131").Value = 8239.764299999999
$ws.Range("M131").Value = -16042233
$ws.Range("N131").Value = -18319.7643

# ===== Sheet 6: GSM =====
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 145.55
$ws.Range("I2").Value = 146.1
$ws.Range("J2").Value = 145
$ws.Range("K2").Value = 146.1
$ws.Range("L2").Value = 145
$ws.Range("M2").Value = -33.09999999999999
$ws.Range("N2").Value = -371
$ws.Range("H6").Value = 1500
$ws.Range("I6").Value = 1000
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -887
$ws.Range("N6").Value = -2226
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -750
$ws.Range("N16").Value = -2500
$ws.Range("H70").Value = 6483.25
$ws.Range("I70").Value = 3957.5
$ws.Range("J70").Value = 9009
$ws.Range("K70").Value = 3957.5
$ws.Range("L70").Value = 9009
$ws.Range("M70").Value = -3687.5
$ws.Range("N70").Value = -9549
$ws.Range("H73").Value = 6483.25
$ws.Range("I73").Value = 3957.5
$ws.Range("J73").Value = 9009
$ws.Range("K73").Value = 3957.5
$ws.Range("L73").Value = 9009
$ws.Range("M73").Value = -3021.5
$ws.Range("N73").Value = -10881
$ws.Range("H80").Value = 21843916
$ws.Range("I80").Value = 126612.78
$ws.Range("J80").Value = 35805040
$ws.Range("K80").Value = 126612.78
$ws.Range("L80").Value = 35805040
$ws.Range("M80").Value = -125614.78
$ws.Range("N80").Value = -35807036
$ws.Range("H83").Value = 21843916
$ws.Range("I83").Value = 126612.78
$ws.Range("J83").Value = 35805040
$ws.Range("K83").Value = 633063.9
$ws.Range("L83").Value = 179025200
$ws.Range("M83").Value = -628071.9
$ws.Range("N83").Value = -179035184
$ws.Range("H97").Value = 6740
$ws.Range("I97").Value = 1201
$ws.Range("J97").Value = 15444.143
$ws.Range("K97").Value = 1201
$ws.Range("L97").Value = 15444.143
$ws.Range("M97").Value = -705
$ws.Range("N97").Value = -16436.143
$ws.Range("H113").Value = 1377.8
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 5576.2856
$ws.Range("I122").Value = 4473.6
$ws.Range("K122").Value = 13420.8
$ws.Range("M122").Value = -10970.8
$ws.Range("H132").Value = 6412712
$ws.Range("I132").Value = 3275
$ws.Range("J132").Value = 25641024
$ws.Range("K132").Value = 9825
$ws.Range("L132").Value = 76923072
$ws.Range("M132").Value = -7295
$ws.Range("N132").Value = -76928132

# ===== Sheet 7: LTW =====
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 12865.762
$ws.Range("I7").Value = 8252.444
$ws.Range("J7").Value = 16325.75
$ws.Range("K7").Value = 8252.444
$ws.Range("L7").Value = 16325.75
$ws.Range("M7").Value = -8140.444
$ws.Range("N7").Value = -16549.75
$ws.Range("H13").Value = 52100.297
$ws.Range("I13").Value = 8699.75
$ws.Range("K13").Value = 8699.75
$ws.Range("M13").Value = -8559.75
$ws.Range("H16").Value = 837.9231
$ws.Range("I16").Value = 735.7273
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 735.7273
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -565.7273
$ws.Range("N16").Value = -1740
$ws.Range("H22").Value = 5437.9
$ws.Range("J22").Value = 6522.625
$ws.Range("L22").Value = 6522.625
$ws.Range("N22").Value = -7112.625
$ws.Range("H27").Value = 5437.9
$ws.Range("J27").Value = 6522.625
$ws.Range("L27").Value = 6522.625
$ws.Range("N27").Value = -6736.625
$ws.Range("H40").Value = 2141432
$ws.Range("I40").Value = 3271291
$ws.Range("J40").Value = 7254.1113
$ws.Range("K40").Value = 3271291
$ws.Range("L40").Value = 7254.1113
$ws.Range("M40").Value = -3271155
$ws.Range("N40").Value = -7526.1113
$ws.Range("H46").Value = 9248.862999999999
$ws.Range("I46").Value = 28373
$ws.Range("K46").Value = 28373
$ws.Range("M46").Value = -28185
$ws.Range("H61").Value = 12504279
$ws.Range("I61").Value = 16671163
$ws.Range("K61").Value = 16671163
$ws.Range("M61").Value = -16670961
$ws.Range("H63").Value = 41453.89
$ws.Range("I63").Value = 39012.145
$ws.Range("K63").Value = 39012.145
$ws.Range("M63").Value = -38263.145
$ws.Range("H66").Value = 41453.89
$ws.Range("I66").Value = 39012.145
$ws.Range("K66").Value = 117036.435
$ws.Range("M66").Value = -113292.435
$ws.Range("H68").Value = 1930.8334
$ws.Range("I68").Value = 1241.2222
$ws.Range("K68").Value = 1241.2222
$ws.Range("M68").Value = -492.2221999999999
$ws.Range("H71").Value = 1930.8334
$ws.Range("I71").Value = 1241.2222
$ws.Range("K71").Value = 6206.111
$ws.Range("M71").Value = -2462.111
$ws.Range("H74").Value = 43944.168
$ws.Range("I74").Value = 41925.555
$ws.Range("K74").Value = 41925.555
$ws.Range("M74").Value = -40927.555
$ws.Range("H77").Value = 43944.168
$ws.Range("I77").Value = 41925.555
$ws.Range("K77").Value = 125776.665
$ws.Range("M77").Value = -120784.665
$ws.Range("H82").Value = 2135.2778
$ws.Range("I82").Value = 2122.6
$ws.Range("K82").Value = 2122.6
$ws.Range("M82").Value = -1761.6
$ws.Range("H85").Value = 2135.2778
$ws.Range("I85").Value = 2122.6
$ws.Range("K85").Value = 2122.6
$ws.Range("M85").Value = -874.5999999999999
$ws.Range("H93").Value = 1285.7894
$ws.Range("I93").Value = 1301.6666
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1301.6666
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -53.66660000000002
$ws.Range("N93").Value = -3496
$ws.Range("H113").Value = 12504279
$ws.Range("I113").Value = 16671163
$ws.Range("K113").Value = 16671163
$ws.Range("M113").Value = -16668993
$ws.Range("H122").Value = 4057.7778
$ws.Range("I122").Value = 3018.1667
$ws.Range("J122").Value = 5245.905
$ws.Range("K122").Value = 9054.500100000001
$ws.Range("L122").Value = 15737.715
$ws.Range("M122").Value = -6604.500100000001
$ws.Range("N122").Value = -20637.715
$ws.Range("H126").Value = 12865.762
$ws.Range("I126").Value = 8252.444
$ws.Range("J126").Value = 16325.75
$ws.Range("K126").Value = 24757.332
$ws.Range("L126").Value = 48977.25
$ws.Range("M126").Value = -22287.332
$ws.Range("N126").Value = -53917.25
$ws.Range("H132").Value = 2215.94
$ws.Range("I132").Value = 2252.8823
$ws.Range("J132").Value = 2137.4375
$ws.Range("K132").Value = 6758.646900000001
$ws.Range("L132").Value = 6412.3125
$ws.Range("M132").Value = -4228.646900000001
$ws.Range("N132").Value = -11472.3125
$ws.Range("H136").Value = 4523.5386
$ws.Range("I136").Value = 2862.875
$ws.Range("J136").Value = 7180.6
$ws.Range("K136").Value = 8588.625
$ws.Range("L136").Value = 21541.8
$ws.Range("M136").Value = -6038.625
$ws.Range("N136").Value = -26641.8

# ===== Sheet 8: WVR =====
$ws = $wb.Worksheets.Item(8)
$ws.Range("H69").Value = 24016.5
$ws.Range("J69").Value = 18819.8
$ws.Range("L69").Value = 18819.8
$ws.Range("N69").Value = -20317.8
$ws.Range("H72").Value = 24016.5
$ws.Range("J72").Value = 18819.8
$ws.Range("L72").Value = 56459.39999999999
$ws.Range("N72").Value = -63947.39999999999
$ws.Range("H96").Value = 7283.8335
$ws.Range("J96").Value = 8501.5
$ws.Range("L96").Value = 8501.5
$ws.Range("N96").Value = -11247.5
$ws.Range("H100").Value = 981020.0600000001
$ws.Range("I100").Value = 1508795.1
$ws.Range("J100").Value = 866.2857
$ws.Range("K100").Value = 3017590.2
$ws.Range("L100").Value = 1732.5714
$ws.Range("M100").Value = -3017049.2
$ws.Range("N100").Value = -2814.5714
$ws.Range("H105").Value = 54903.5
$ws.Range("J105").Value = 54903.5
$ws.Range("L105").Value = 54903.5
$ws.Range("N105").Value = -61891.5
$ws.Range("H113").Value = 428.2857
$ws.Range("I113").Value = 373.75
$ws.Range("J113").Value = 501
$ws.Range("K113").Value = 1121.25
$ws.Range("L113").Value = 1503
$ws.Range("M113").Value = 1048.75
$ws.Range("N113").Value = -5843
$ws.Range("H122").Value = 2287.1155
$ws.Range("I122").Value = 2200.7144
$ws.Range("J122").Value = 2650
$ws.Range("K122").Value = 6602.1432
$ws.Range("L122").Value = 7950
$ws.Range("M122").Value = -4152.1432
$ws.Range("N122").Value = -12850
$ws.Range("H132").Value = 33199
$ws.Range("I132").Value = 57309.445
$ws.Range("J132").Value = 2199.8572
$ws.Range("K132").Value = 171928.335
$ws.Range("L132").Value = 6599.571599999999
$ws.Range("M132").Value = -169398.335
$ws.Range("N132").Value = -11659.5716
$ws.Range("H135").Value = 92330.164
$ws.Range("J135").Value = 92330.164
$ws.Range("L135").Value = 92330.164
$ws.Range("N135").Value = -102470.164
$ws.Range("H136").Value = 50148
$ws.Range("I136").Value = 65209.94
$ws.Range("K136").Value = 195629.82
$ws.Range("M136").Value = -193079.82
